# This script reproduces the crypto price/volume update described in the commit.
# Cells in column D that are plain decimal numbers must be forced to remain text
# (matching the original inline-string cell type), otherwise Excel auto-converts
# them to numeric values and changes their displayed precision (e.g. trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.951.38"
$ws.Range("E2").Value = "  +6.45%  "

$ws.Range("D3").Value = "1.738.15"
$ws.Range("E3").Value = "  +4.95%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'229.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.64%  "

$ws.Range("D6").Value = "'0.5444"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.85%  "

$ws.Range("D7").Value = "'1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.21%  "

$ws.Range("D8").Value = "'0.2781"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.20%  "

$ws.Range("D9").Value = "'0.06713"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.46%  "

$ws.Range("D10").Value = "'21.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.00%  "

$ws.Range("D11").Value = "'0.07780"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("D12").Value = "'4.713"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.44%  "

$ws.Range("D13").Value = "1.740.11"
$ws.Range("E13").Value = "  +4.43%  "

$ws.Range("D14").Value = "1.977.26"
$ws.Range("E14").Value = "  +4.90%  "

$ws.Range("D15").Value = "'0.5998"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.31%  "

$ws.Range("D16").Value = "0.0₅8427"
$ws.Range("E16").Value = "  +2.00%  "

$ws.Range("D17").Value = "'69.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.44%  "

$ws.Range("D18").Value = "27.935.52"
$ws.Range("E18").Value = "  +6.36%  "

$ws.Range("D19").Value = "'224.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +16.40%  "

$ws.Range("D20").Value = "'4.833"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.91%  "

$ws.Range("D21").Value = "'1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").Value = "'10.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.77%  "

$ws.Range("E23").Value = "  +4.08%  "

$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").Value = "'146.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.25%  "

$ws.Range("E26").Value = "  +3.64%  "

$ws.Range("D27").Value = "'7.461"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.42%  "

$ws.Range("E28").Value = "  +7.41%  "

$ws.Range("D29").Value = "'1.649"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.04%  "

$ws.Range("D30").Value = "'0.05654"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.16%  "

$ws.Range("D31").Value = "'1.320"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.96%  "

$ws.Range("D32").Value = "'3.699"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.35%  "

$ws.Range("D33").Value = "'3.549"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.43%  "

$ws.Range("D34").Value = "'1.663"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.87%  "

$ws.Range("D35").Value = "'0.9856"
$ws.Range("D35").Style = "Normal"

$ws.Range("D37").Value = "'2.451"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.53%  "

$ws.Range("D38").Value = "'0.5953"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.02%  "

$ws.Range("D39").Value = "'0.01675"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.53%  "

$ws.Range("D40").Value = "'6.027"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.68%  "

$ws.Range("D41").Value = "'0.8480"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("D42").Value = "1.047.83"
$ws.Range("E42").Value = "  +3.16%  "

$ws.Range("D43").Value = "'1.004"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("D44").Value = "'102.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.39%  "

$ws.Range("D45").Value = "1.883.53"
$ws.Range("E45").Value = "  +4.85%  "

$ws.Range("E46").Value = "  +3.09%  "

$ws.Range("E47").Value = "  -1.66%  "

$ws.Range("D48").Value = "'8.321"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.65%  "

$ws.Range("D49").Value = "'1.016"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").Value = "'0.4428"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.76%  "

$ws.Range("D51").Value = "'0.05318"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.63%  "
